# Insert a new data row at row 375 (pushes the old rows 375-446 down to 376-447,
# growing the used range from A1:R446 to A1:R447), then populate the new row
# with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 375, shifting everything below it down by one.
$ws.Rows(375).Insert()

# Fill in the new row 375 with its data.
$ws.Cells.Item(375, 1).Value = 4
$ws.Cells.Item(375, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(375, 3).Value = "Los Lagos"
$ws.Cells.Item(375, 4).Value = 44711
$ws.Cells.Item(375, 5).Value = 10
$ws.Cells.Item(375, 6).Value = 100112006
$ws.Cells.Item(375, 7).Value = "Repollo"
$ws.Cells.Item(375, 8).Value = "Crespo record"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 750
$ws.Cells.Item(375, 11).Value = 1500
$ws.Cells.Item(375, 12).Value = 1500
$ws.Cells.Item(375, 13).Value = 1500
$ws.Cells.Item(375, 14).Value = "$/unidad"
$ws.Cells.Item(375, 15).Value = "Región del Maule"
$ws.Cells.Item(375, 16).Value = 1500
$ws.Cells.Item(375, 17).Value = 1
$ws.Cells.Item(375, 18).Value = "Hortaliza"
